$wb = $excel.ActiveWorkbook

# Rename the "IAir" sheet to "Air"
$ws = $wb.Worksheets.Item("IAir")
$ws.Name = "Air"

# The "s^0" header (shared with an "x" quality-column header elsewhere in the
# workbook) becomes "s0" -- setting the header cell's value updates the
# shared-string table (and, since the cell anchors the table's 6th column,
# the table's column name) in one shot.
$ws.Range("F2").Value = "s0"

# Make "Air" the active sheet/tab and select H33 on it.
$ws.Activate()
$ws.Range("H33").Select()
